# Update header labels in row 1 (columns E:L) on every worksheet to use
# the more descriptive "severity level(s)" wording.

$wb = $excel.ActiveWorkbook

$replacements = @{
    "% 1-2" = "% severity levels 1-2"
    "# 1-2" = "# severity levels 1-2"
    "% 3"   = "% severity level 3"
    "# 3"   = "# severity level 3"
    "% 4"   = "% severity level 4"
    "# 4"   = "# severity level 4"
    "% 5"   = "% severity level 5"
    "# 5"   = "# severity level 5"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($col in "E", "F", "G", "H", "I", "J", "K", "L") {
        $cell = $ws.Range($col + "1")
        $current = $cell.Value2
        if ($replacements.ContainsKey($current)) {
            $cell.Value = $replacements[$current]
        }
    }
}
